$d = $word.ActiveDocument

# --- Update the date heading paragraph ---
$d.Content.Find.Execute("2025-08-31 Sunday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-09-01 Monday", 2)

# --- Update the division problems in the table ---
# The table has 20 rows x 5 columns; data lives in rows 1, 5, 9, 13, 17.
# Each row maps to an array of [oldText -> newText] pairs, applied by
# direct cell addressing so no text collisions/cascades can occur.

$t = $d.Tables.Item(1)

# Directly replace the run text by setting Range.Text for the cell range
# (excluding the end-of-cell marker), addressed by row/column. This avoids
# any find/replace text-collision issues between cells.

function Set-Cell($row, $col, $text) {
    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 1   # exclude end-of-cell marker
    $r.Text = $text
}

Set-Cell 1 1 "174÷9="
Set-Cell 1 2 "183÷5="
Set-Cell 1 3 "522÷9="
Set-Cell 1 4 "797÷7="
Set-Cell 1 5 "212÷8="

Set-Cell 5 1 "305÷4="
Set-Cell 5 2 "581÷2="
Set-Cell 5 3 "683÷2="
Set-Cell 5 4 "153÷8="
Set-Cell 5 5 "509÷9="

Set-Cell 9 1 "335÷2="
Set-Cell 9 2 "697÷4="
Set-Cell 9 3 "118÷2="
Set-Cell 9 4 "658÷5="
Set-Cell 9 5 "154÷5="

Set-Cell 13 1 "478÷6="
Set-Cell 13 2 "123÷4="
Set-Cell 13 3 "701÷7="
Set-Cell 13 4 "415÷7="
Set-Cell 13 5 "664÷9="

Set-Cell 17 1 "848÷8="
Set-Cell 17 2 "458÷8="
Set-Cell 17 3 "212÷2="
Set-Cell 17 4 "631÷2="
Set-Cell 17 5 "456÷2="
